$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formats from row 18 down to row 19 so the new row matches existing styling
$ws.Range("A18:G18").Copy() | Out-Null
$ws.Range("A19:G19").PasteSpecial(-4122) | Out-Null

# Add new row 19: a new timekeeping entry
$ws.Range("A19").Value = 44173
$ws.Range("B19").Value = 0.56180555555555556
$ws.Range("C19").Value = 0.65486111111111112
$ws.Range("D19").Value = 0

# Extend the shared formula from E2:E18 down to E19
$ws.Range("E19").FormulaR1C1 = "=RC[-2]-RC[-3]-TIME(0,RC[-1],0)"

# New comment string for the added row
$ws.Range("F19").Value = "creating final UML diagrams and report"

# Move the active selection to F20
$ws.Range("F20").Select()
